# [Shreya] Add: Added Flipkart Process and refactored Business object
# Add a "Colour" column (D) to the product listing sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Colour"

# Leave the active selection on the newly added header cell's data start,
# matching the authored workbook's saved selection state.
$ws.Range("D2").Select()
